# A new daily price record was added to the "Pomelo" sheet. In the
# canonical OOXML this shows up as a brand-new row 301 (Fecha=44985,
# Volumen=55) with every following row (old 301..388) shifted down by
# one (new 302..389), and the sheet's used range growing from
# A1:T388 to A1:T389.
#
# Reproduce that with a real row insert (so everything below shifts
# down automatically) followed by populating the newly blank row with
# the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 301 - this pushes old rows 301..388 down
# to 302..389, exactly matching the diff.
$ws.Rows.Item(301).EntireRow.Insert()

# Populate the new row 301 with the new record. Columns A,B,C,E..L,N..T
# mirror the existing series (same market/region/product/variety/
# quality/price-tier/unit/origin); only Fecha (D) and Volumen (M) are
# genuinely new values per the diff.
$ws.Cells.Item(301, 1).Value = 10
$ws.Cells.Item(301, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(301, 3).Value = "La Araucanía"
$ws.Cells.Item(301, 4).Value = 44985
$ws.Cells.Item(301, 5).Value = 9
$ws.Cells.Item(301, 6).Value = "Fruta"
$ws.Cells.Item(301, 7).Value = 100102
$ws.Cells.Item(301, 8).Value = "Cítricos"
$ws.Cells.Item(301, 9).Value = 100102006
$ws.Cells.Item(301, 10).Value = "Pomelo"
$ws.Cells.Item(301, 11).Value = "Start Ruby"
$ws.Cells.Item(301, 12).Value = "Primera"
$ws.Cells.Item(301, 13).Value = 55
$ws.Cells.Item(301, 14).Value = 14000
$ws.Cells.Item(301, 15).Value = 14000
$ws.Cells.Item(301, 16).Value = 14000
$ws.Cells.Item(301, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(301, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(301, 19).Value = 933
$ws.Cells.Item(301, 20).Value = 15
